$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing benchmark values for row 20 (F20 and G20)
$ws.Range("F20").Value = 1036
$ws.Range("G20").Value = 45595

# Update the active selection to match the edit (F21)
$ws.Range("F21").Select()
